$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 437.9524
$ws.Range("J17").Value = 437.9524
$ws.Range("L17").Value = 1313.8572
$ws.Range("N17").Value = -1649.8572

$ws.Range("H51").Value = 10421602
$ws.Range("J51").Value = 7446.7144
$ws.Range("L51").Value = 7446.7144
$ws.Range("N51").Value = -8414.714400000001

$ws.Range("H62").Value = 5225
$ws.Range("I62").Value = 5250
$ws.Range("J62").Value = 5200
$ws.Range("K62").Value = 5250
$ws.Range("L62").Value = 5200
$ws.Range("M62").Value = -4626
$ws.Range("N62").Value = -6448

$ws.Range("H65").Value = 5225
$ws.Range("I65").Value = 5250
$ws.Range("J65").Value = 5200
$ws.Range("K65").Value = 26250
$ws.Range("L65").Value = 26000
$ws.Range("M65").Value = -23130
$ws.Range("N65").Value = -32240

$ws.Range("H98").Value = 2840.3635
$ws.Range("I98").Value = 2022.9524
$ws.Range("J98").Value = 20006
$ws.Range("K98").Value = 2022.9524
$ws.Range("L98").Value = 20006
$ws.Range("M98").Value = -524.9523999999999
$ws.Range("N98").Value = -23002

$ws.Range("H115").Value = 454.55554
$ws.Range("I115").Value = 417.625
$ws.Range("K115").Value = 1252.875
$ws.Range("M115").Value = 314.125

$ws.Range("H122").Value = 2840.3635
$ws.Range("I122").Value = 2022.9524
$ws.Range("J122").Value = 20006
$ws.Range("K122").Value = 6068.857199999999
$ws.Range("L122").Value = 60018
$ws.Range("M122").Value = -3618.857199999999
$ws.Range("N122").Value = -64918

$ws.Range("H129").Value = 1421.091
$ws.Range("I129").Value = 1008.55554
$ws.Range("J129").Value = 3277.5
$ws.Range("K129").Value = 3025.66662
$ws.Range("L129").Value = 9832.5
$ws.Range("M129").Value = 1974.33338
$ws.Range("N129").Value = -19832.5

$ws.Range("H130").Value = 73250
$ws.Range("J130").Value = 106500
$ws.Range("L130").Value = 106500
$ws.Range("N130").Value = -116540

$ws.Range("H135").Value = 3339.3125
$ws.Range("I135").Value = 3193.5454
$ws.Range("K135").Value = 28741.9086
$ws.Range("M135").Value = -26206.9086

$ws.Range("H137").Value = 19125.5
$ws.Range("I137").Value = 2003
$ws.Range("J137").Value = 27686.75
$ws.Range("K137").Value = 6009
$ws.Range("L137").Value = 83060.25
$ws.Range("M137").Value = -3459
$ws.Range("N137").Value = -88160.25

$ws.Range("H138").Value = 1815167.1
$ws.Range("J138").Value = 2690100.2
$ws.Range("L138").Value = 8070300.600000001
$ws.Range("N138").Value = -8080580.600000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15061.13
$ws.Range("I32").Value = 15318.567
$ws.Range("K32").Value = 15318.567
$ws.Range("M32").Value = -15031.567

$ws.Range("H122").Value = 1947.3478
$ws.Range("I122").Value = 1637.8125
$ws.Range("J122").Value = 2654.8572
$ws.Range("K122").Value = 4913.4375
$ws.Range("L122").Value = 7964.571599999999
$ws.Range("M122").Value = -2463.4375
$ws.Range("N122").Value = -12864.5716

$ws.Range("H132").Value = 1813.9744
$ws.Range("I132").Value = 1461.7576
$ws.Range("J132").Value = 3751.1667
$ws.Range("K132").Value = 4385.2728
$ws.Range("L132").Value = 11253.5001
$ws.Range("M132").Value = -1855.2728
$ws.Range("N132").Value = -16313.5001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 371.66666
$ws.Range("I12").Value = 371.66666
$ws.Range("K12").Value = 371.66666
$ws.Range("M12").Value = -203.66666

$ws.Range("H82").Value = 13502.8
$ws.Range("I82").Value = 8128.5
$ws.Range("K82").Value = 8128.5
$ws.Range("M82").Value = -7745.5

$ws.Range("H85").Value = 13502.8
$ws.Range("I85").Value = 8128.5
$ws.Range("K85").Value = 8128.5
$ws.Range("M85").Value = -6802.5

$ws.Range("H134").Value = 3372.6
$ws.Range("I134").Value = 3253.1333
$ws.Range("K134").Value = 9759.3999
$ws.Range("M134").Value = -7224.3999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1820516.4
$ws.Range("I31").Value = 2858947.8
$ws.Range("J31").Value = 3261.2
$ws.Range("K31").Value = 2858947.8
$ws.Range("L31").Value = 3261.2
$ws.Range("M31").Value = -2858652.8
$ws.Range("N31").Value = -3851.2

$ws.Range("H34").Value = 1820516.4
$ws.Range("I34").Value = 2858947.8
$ws.Range("J34").Value = 3261.2
$ws.Range("K34").Value = 2858947.8
$ws.Range("L34").Value = 3261.2
$ws.Range("M34").Value = -2858745.8
$ws.Range("N34").Value = -3665.2

$ws.Range("H58").Value = 1464.5385
$ws.Range("J58").Value = 1721.2941
$ws.Range("L58").Value = 1721.2941
$ws.Range("N58").Value = -2127.2941

$ws.Range("H132").Value = 3479.4546
$ws.Range("I132").Value = 2314.8333
$ws.Range("J132").Value = 4877
$ws.Range("K132").Value = 6944.499899999999
$ws.Range("L132").Value = 14631
$ws.Range("M132").Value = -4414.499899999999
$ws.Range("N132").Value = -19691

$ws.Range("H136").Value = 1464.5385
$ws.Range("J136").Value = 1721.2941
$ws.Range("L136").Value = 5163.8823
$ws.Range("N136").Value = -10263.8823

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3249.157
$ws.Range("I68").Value = 1668.4286
$ws.Range("J68").Value = 3500.6365
$ws.Range("K68").Value = 5005.2858
$ws.Range("L68").Value = 10501.9095
$ws.Range("M68").Value = -4194.2858
$ws.Range("N68").Value = -12123.9095

$ws.Range("H71").Value = 3249.157
$ws.Range("I71").Value = 1668.4286
$ws.Range("J71").Value = 3500.6365
$ws.Range("K71").Value = 15015.8574
$ws.Range("L71").Value = 31505.7285
$ws.Range("M71").Value = -10959.8574
$ws.Range("N71").Value = -39617.7285

$ws.Range("H113").Value = 785.2308
$ws.Range("J113").Value = 817.9167
$ws.Range("L113").Value = 2453.7501
$ws.Range("N113").Value = -6793.7501

$ws.Range("H129").Value = 3937.8333
$ws.Range("I129").Value = 3627.5
$ws.Range("K129").Value = 10882.5
$ws.Range("M129").Value = -5882.5

$ws.Range("H131").Value = 2553.432
$ws.Range("I131").Value = 4064.4546
$ws.Range("J131").Value = 2049.7576
$ws.Range("K131").Value = 12193.3638
$ws.Range("L131").Value = 6149.2728
$ws.Range("M131").Value = -7153.363799999999
$ws.Range("N131").Value = -16229.2728

$ws.Range("H137").Value = 1855.1666
$ws.Range("J137").Value = 3016
$ws.Range("L137").Value = 9048
$ws.Range("N137").Value = -19248

$ws.Range("H139").Value = 3190.5715
$ws.Range("I139").Value = 2111.2778
$ws.Range("K139").Value = 6333.8334
$ws.Range("M139").Value = -1193.8334

$ws.Range("H140").Value = 84829
$ws.Range("I140").Value = 84829
$ws.Range("K140").Value = 254487
$ws.Range("M140").Value = -249307

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2177.3635
$ws.Range("J102").Value = 2672.25
$ws.Range("L102").Value = 2672.25
$ws.Range("N102").Value = -5916.25

$ws.Range("H122").Value = 38463496
$ws.Range("I122").Value = 1116
$ws.Range("K122").Value = 3348
$ws.Range("M122").Value = -898

$ws.Range("H126").Value = 1704.4546
$ws.Range("I126").Value = 958.3333
$ws.Range("J126").Value = 2599.8
$ws.Range("K126").Value = 2874.9999
$ws.Range("L126").Value = 7799.400000000001
$ws.Range("M126").Value = -404.9998999999998
$ws.Range("N126").Value = -12739.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2662.3428
$ws.Range("I40").Value = 2611.7812
$ws.Range("J40").Value = 3201.6667
$ws.Range("K40").Value = 2611.7812
$ws.Range("L40").Value = 3201.6667
$ws.Range("M40").Value = -2475.7812
$ws.Range("N40").Value = -3473.6667

$ws.Range("H82").Value = 3634.6843
$ws.Range("J82").Value = 4851.857
$ws.Range("L82").Value = 4851.857
$ws.Range("N82").Value = -5573.857

$ws.Range("H85").Value = 3634.6843
$ws.Range("J85").Value = 4851.857
$ws.Range("L85").Value = 4851.857
$ws.Range("N85").Value = -7347.857

$ws.Range("H94").Value = 45000
$ws.Range("J94").Value = 45000
$ws.Range("L94").Value = 45000
$ws.Range("N94").Value = -46352

$ws.Range("H122").Value = 4997.25
$ws.Range("I122").Value = 4997
$ws.Range("K122").Value = 14991
$ws.Range("M122").Value = -12541

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 19397.4
$ws.Range("I62").Value = 20000
$ws.Range("K62").Value = 20000
$ws.Range("M62").Value = -19376

$ws.Range("H65").Value = 19397.4
$ws.Range("I65").Value = 20000
$ws.Range("K65").Value = 100000
$ws.Range("M65").Value = -96880

$ws.Range("H114").Value = 72000
$ws.Range("J114").Value = 72000
$ws.Range("L114").Value = 72000
$ws.Range("N114").Value = -80678

$ws.Range("H115").Value = 32222.223
$ws.Range("J115").Value = 32222.223
$ws.Range("L115").Value = 32222.223
$ws.Range("N115").Value = -35356.223

$ws.Range("H116").Value = 90500
$ws.Range("J116").Value = 90500
$ws.Range("L116").Value = 90500
$ws.Range("N116").Value = -99678

$ws.Range("H122").Value = 4334.485
$ws.Range("I122").Value = 4939.524
$ws.Range("J122").Value = 3275.6667
$ws.Range("K122").Value = 14818.572
$ws.Range("L122").Value = 9827.000100000001
$ws.Range("M122").Value = -12368.572
$ws.Range("N122").Value = -14727.0001

$ws.Range("H126").Value = 2603.7778
$ws.Range("I126").Value = 1956.6
$ws.Range("J126").Value = 4452.857
$ws.Range("K126").Value = 5869.799999999999
$ws.Range("L126").Value = 13358.571
$ws.Range("M126").Value = -3399.799999999999
$ws.Range("N126").Value = -18298.571

$ws.Range("H132").Value = 26849.107
$ws.Range("I132").Value = 30848.666
$ws.Range("J132").Value = 2851.75
$ws.Range("K132").Value = 92545.99800000001
$ws.Range("L132").Value = 8555.25
$ws.Range("M132").Value = -90015.99800000001
$ws.Range("N132").Value = -13615.25
